{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the \"Added many more features\" commit:\n//   - Updates the title (heading + bold CTA run) to drop \"for\" and use a\n//     plain hyphen instead of an en dash.\n//   - Rewrites the \"What we like\" bullet list with new copy.\n//   - Rewrites two \"What we don't like\" bullets with new copy.\n//   - Shortens the italic meta-description paragraph.\n\n// Each entry is an exact, literal (case-sensitive) run of text found via\n// Word's search API and replaced in place, so paragraph styles and run\n// formatting (bold/italic/etc.) on the surrounding paragraph are preserved.\nconst replacements = [\n  {\n    find: \"Play Bigger Bass Blizzard \\u2013 Christmas Catch for Free\",\n    replace: \"Play Bigger Bass Blizzard - Christmas Catch Free\",\n  },\n  {\n    find: \"Stunning graphics and Christmas-themed music\",\n    replace: \"Stunning graphics and visual design\",\n  },\n  {\n    find: \"High maximum win of 2,667 times the bet\",\n    replace: \"Festive Christmas-themed background music\",\n  },\n  {\n    find: \"Turbo and Autoplay features for faster play\",\n    replace: \"High volatility for exciting gameplay\",\n  },\n  {\n    find: \"Free Spins feature with a Wild symbol and multiplier\",\n    replace: \"Chance to win up to 4,000 times the bet\",\n  },\n  {\n    find: \"Only 12 fixed paylines with low stakes\",\n    replace: \"Limited number of paylines\",\n  },\n  {\n    find: \"High volatility may not suit inexperienced players\",\n    replace: \"High volatility may not appeal to all players\",\n  },\n  {\n    find:\n      \"Read our review of Bigger Bass Blizzard \\u2013 Christmas Catch and play for free. Discover its stunning graphics, high maximum win, and special Free Spins feature.\",\n    replace:\n      \"Read our review of Bigger Bass Blizzard - Christmas Catch and play for free.\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the \"Added many more features\" commit:\n#   - Updates the title (heading + bold CTA run) to drop \"for\" and use a\n#     plain hyphen instead of an en dash.\n#   - Rewrites the \"What we like\" bullet list with new copy.\n#   - Rewrites two \"What we don't like\" bullets with new copy.\n#   - Shortens the italic meta-description paragraph.\n#\n# Each pair is an exact, literal (case-sensitive) string replaced via\n# Range.Find/Execute with ReplaceAll so paragraph styles and surrounding run\n# formatting (bold/italic/etc.) are left untouched.\n\n$d = $word.ActiveDocument\n\n$enDash = [char]0x2013\n\n$replacements = @(\n    @{\n        Find    = \"Play Bigger Bass Blizzard $enDash Christmas Catch for Free\"\n        Replace = \"Play Bigger Bass Blizzard - Christmas Catch Free\"\n    },\n    @{\n        Find    = \"Stunning graphics and Christmas-themed music\"\n        Replace = \"Stunning graphics and visual design\"\n    },\n    @{\n        Find    = \"High maximum win of 2,667 times the bet\"\n        Replace = \"Festive Christmas-themed background music\"\n    },\n    @{\n        Find    = \"Turbo and Autoplay features for faster play\"\n        Replace = \"High volatility for exciting gameplay\"\n    },\n    @{\n        Find    = \"Free Spins feature with a Wild symbol and multiplier\"\n        Replace = \"Chance to win up to 4,000 times the bet\"\n    },\n    @{\n        Find    = \"Only 12 fixed paylines with low stakes\"\n        Replace = \"Limited number of paylines\"\n    },\n    @{\n        Find    = \"High volatility may not suit inexperienced players\"\n        Replace = \"High volatility may not appeal to all players\"\n    },\n    @{\n        Find    = \"Read our review of Bigger Bass Blizzard $enDash Christmas Catch and play for free. Discover its stunning graphics, high maximum win, and special Free Spins feature.\"\n        Replace = \"Read our review of Bigger Bass Blizzard - Christmas Catch and play for free.\"\n    }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Execute($r.Find, $true, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2) | Out-Null\n}\n"}
